$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E value cells are written as Text so that numeric-looking
# strings (e.g. "1.00", "0.140") keep their exact original formatting
# instead of being auto-converted to numbers by Excel.
$valueRange = $ws.Range("D2:E51")
$valueRange.NumberFormat = "@"

$ws.Range("D2").Value = "96.875.59"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "3.317.75"
$ws.Range("E3").Value = "  -4.24%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "246.53"
$ws.Range("E5").Value = "  -6.03%  "
$ws.Range("D6").Value = "650.58"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("D7").Value = "1.34"
$ws.Range("E7").Value = "  -14.68%  "
$ws.Range("D8").Value = "0.409"
$ws.Range("E8").Value = "  -10.87%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "0.968"
$ws.Range("E10").Value = "  -14.80%  "
$ws.Range("D11").Value = "3.317.40"
$ws.Range("E11").Value = "  -4.15%  "
$ws.Range("D12").Value = "0.203"
$ws.Range("E12").Value = "  -7.38%  "
$ws.Range("D13").Value = "39.53"
$ws.Range("E13").Value = "  -7.67%  "
$ws.Range("D14").Value = "96.510.00"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "5.94"
$ws.Range("E15").Value = "  -4.09%  "
$ws.Range("D16").Value = "0.0000249"
$ws.Range("E16").Value = "  -9.04%  "
$ws.Range("D17").Value = "3.928.22"
$ws.Range("E17").Value = "  -4.38%  "
$ws.Range("D18").Value = "8.53"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "3.314.34"
$ws.Range("E19").Value = "  -4.41%  "
$ws.Range("D20").Value = "16.63"
$ws.Range("E20").Value = "  -5.59%  "
$ws.Range("B21").Value = "SuiNetwork"
$ws.Range("C21").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D21").Value = "3.34"
$ws.Range("E21").Value = "  -7.89%  "
$ws.Range("D22").Value = "10.32"
$ws.Range("E22").Value = "  -4.76%  "
$ws.Range("D23").Value = "491.37"
$ws.Range("E23").Value = "  -7.49%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").Value = "0.458"
$ws.Range("E24").Value = "  -5.49%  "
$ws.Range("E25").Value = "  -10.22%  "
$ws.Range("D26").Value = "6.37"
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("D27").Value = "92.03"
$ws.Range("E27").Value = "  -10.79%  "
$ws.Range("E28").Value = "  -7.93%  "
$ws.Range("D29").Value = "3.500.12"
$ws.Range("E29").Value = "  -4.27%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "0.140"
$ws.Range("E31").Value = "  -8.69%  "
$ws.Range("D32").Value = "10.68"
$ws.Range("E32").Value = "  -8.28%  "
$ws.Range("E33").Value = "  -7.13%  "
$ws.Range("D34").Value = "2.43"
$ws.Range("E34").Value = "  +8.70%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "0.538"
$ws.Range("E36").Value = "  -6.44%  "
$ws.Range("D37").Value = "27.72"
$ws.Range("E37").Value = "  -9.19%  "
$ws.Range("E38").Value = "  +2.49%  "
$ws.Range("D39").Value = "7.45"
$ws.Range("E39").Value = "  -7.37%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "0.148"
$ws.Range("E41").Value = "  -8.05%  "
$ws.Range("D42").Value = "497.90"
$ws.Range("E42").Value = "  -7.41%  "
$ws.Range("D43").Value = "24.49"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").Value = "3.65"
$ws.Range("E44").Value = "  -2.67%  "
$ws.Range("D45").Value = "0.815"
$ws.Range("E45").Value = "  -5.87%  "
$ws.Range("D46").Value = "0.0402"
$ws.Range("E46").Value = "  -8.68%  "
$ws.Range("D47").Value = "8.26"
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("D48").Value = "1.61"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").Value = "5.35"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").Value = "52.39"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").Value = "3.09"
$ws.Range("E51").Value = "  -11.78%  "

# Restore cells to their original (unformatted) style so the
# cells keep behaving like plain/general string cells.
$valueRange.ClearFormats()
